$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 27862
$ws.Range("B2").Value = "Luiz Fernando Machado"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45092
$ws.Range("G2").Value = 3031.32

# Row 3
$ws.Range("A3").Value = 35798
$ws.Range("B3").Value = "Benjamim Machado"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45104
$ws.Range("G3").Value = 2345.18

# Row 4
$ws.Range("A4").Value = 83132
$ws.Range("B4").Value = "Manuela Rocha"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45092
$ws.Range("G4").Value = 9835.440000000001

# Row 5
$ws.Range("A5").Value = 83783
$ws.Range("B5").Value = "Apollo Costa"
$ws.Range("C5").Value = "P&D"
$ws.Range("F5").Value = 45094
$ws.Range("G5").Value = 4086.58

# Row 6
$ws.Range("A6").Value = 16656
$ws.Range("B6").Value = "Dr. Bryan Cirino"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 8813.870000000001

# Row 7
$ws.Range("A7").Value = 16826
$ws.Range("B7").Value = "Beatriz Almeida"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Consulta medica"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 45084
$ws.Range("G7").Value = 5391.5

# Row 8
$ws.Range("A8").Value = 24515
$ws.Range("B8").Value = "Ravy Santos"
$ws.Range("C8").Value = "Operacoes"
$ws.Range("D8").Value = "Consulta medica"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45103
$ws.Range("G8").Value = 5765.86

# Row 9
$ws.Range("A9").Value = 13636
$ws.Range("B9").Value = "Dr. Noah Pastor"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 45091
$ws.Range("G9").Value = 3595.17

# Row 10
$ws.Range("A10").Value = 56176
$ws.Range("B10").Value = "Camila Gomes"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45097
$ws.Range("G10").Value = 3621.74

# Row 11
$ws.Range("A11").Value = 1908
$ws.Range("B11").Value = "Dra. Gabrielly Moura"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Viagem de negocios"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45091
$ws.Range("G11").Value = 3055.11
